# 11.7.1.1a.xlsx — add the "2020" year column (I) to the data table,
# mirroring the formatting already used by the 2019 column (H) /
# the rest of the numeric data columns, and reproduce the refreshed
# selection left behind by the editor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Header cell I4 = 2020 -> same look as H4 (D4:H4 "year" header style)
# ---------------------------------------------------------------
$ws.Range("I4").Value = 2020
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------
# I5 ("Кыргызская Республика" row) = 1.5 -> same look as H5, plus
# the "0.0" one-decimal number format the new column introduces.
# ---------------------------------------------------------------
$ws.Range("I5").Value = 1.5
$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I5").NumberFormat = "0.0"

# ---------------------------------------------------------------
# I6:I13 (oblast rows) -> same look as the other plain data cells
# in that style family (e.g. H13), plus the "0.0" number format.
# ---------------------------------------------------------------
$oblastValues = @{
    6  = 0.2
    7  = 0.8
    8  = 0.4
    9  = 1.8
    10 = 0.5
    11 = 0.7
    12 = 1.9
    13 = 4.5999999999999996
}

foreach ($row in $oblastValues.Keys) {
    $cell = "I$row"
    $ws.Range($cell).Value = $oblastValues[$row]
    $ws.Range("H13").Copy()
    $ws.Range($cell).PasteSpecial(-4122)
    $ws.Range($cell).NumberFormat = "0.0"
}

# ---------------------------------------------------------------
# I14 ("г. Ош" total row) = 0.4 -> same border/font as H14 (bottom
# rule), but vertical-center-only alignment (matches the other
# data rows) and the new "0.0" number format.
# ---------------------------------------------------------------
$ws.Range("I14").Value = 0.4
$ws.Range("H14").Copy()
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("I14").NumberFormat = "0.0"
$ws.Range("I14").HorizontalAlignment = 1   # xlGeneral
$ws.Range("I14").WrapText = $false

# ---------------------------------------------------------------
# Leave the selection where the editor left it.
# ---------------------------------------------------------------
$ws.Range("M9").Select()
